# Applies the "Updated cryptos list ... with GitHub Actions" refresh:
# every Coin row's Price (D) and Volume(1h) (E) text is refreshed to the
# latest scrape, and two coin pairs (Maker/VeChain at rows 42-43, and
# USDe/Arweave at rows 46-47) swapped ranking order, bringing new
# Coin/Link/Price/Volume values with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> { column letter -> new literal text } for every touched cell.
$updates = @{
    2  = @{ "D" = "59.390.00";  "E" = "  +3.50%  " }
    3  = @{ "D" = "3.002.02";   "E" = "  +2.26%  " }
    4  = @{                     "E" = "  +0.06%  " }
    5  = @{ "D" = "563.11";     "E" = "  +2.05%  " }
    6  = @{ "D" = "138.26";     "E" = "  +5.98%  " }
    7  = @{                     "E" = "  -0.23%  " }
    8  = @{                     "E" = "  +2.06%  " }
    9  = @{ "D" = "2.987.15";   "E" = "  +2.16%  " }
    10 = @{                     "E" = "  +3.89%  " }
    11 = @{ "D" = "5.14";       "E" = "  +7.70%  " }
    12 = @{                     "E" = "  +2.80%  " }
    13 = @{                     "E" = "  +4.57%  " }
    14 = @{ "D" = "33.74";      "E" = "  +3.88%  " }
    15 = @{                     "E" = "  +2.39%  " }
    16 = @{ "D" = "3.494.96";   "E" = "  +2.22%  " }
    17 = @{ "D" = "7.25";       "E" = "  +7.68%  " }
    18 = @{ "D" = "2.995.99";   "E" = "  +2.23%  " }
    19 = @{ "D" = "59.341.26";  "E" = "  +3.28%  " }
    20 = @{ "D" = "430.40";     "E" = "  +3.75%  " }
    21 = @{ "D" = "13.66";      "E" = "  +4.94%  " }
    22 = @{                     "E" = "  +6.25%  " }
    23 = @{                     "E" = "  +3.00%  " }
    24 = @{ "D" = "13.34";      "E" = "  +3.84%  " }
    25 = @{ "D" = "80.87";      "E" = "  +2.73%  " }
    26 = @{                     "E" = "  +0.01%  " }
    27 = @{                     "E" = "  +0.23%  " }
    28 = @{                     "E" = "  +9.77%  " }
    29 = @{                     "E" = "  +3.04%  " }
    30 = @{ "D" = "7.78";       "E" = "  +4.50%  " }
    31 = @{                     "E" = "  +3.18%  " }
    32 = @{                     "E" = "  +0.26%  " }
    33 = @{ "D" = "0.0991";     "E" = "  -3.51%  " }
    34 = @{ "D" = "5.94";       "E" = "  +6.17%  " }
    35 = @{ "D" = "0.990";      "E" = "  +6.74%  " }
    36 = @{                     "E" = "  +13.35%  " }
    37 = @{                     "E" = "  +0.05%  " }
    38 = @{                     "E" = "  +1.23%  " }
    39 = @{                     "E" = "  +3.56%  " }
    40 = @{                     "E" = "  +6.75%  " }
    41 = @{ "D" = "401.34";     "E" = "  +6.78%  " }
    42 = @{ "B" = "VeChain"; "C" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; "D" = "0.0352";    "E" = "  +1.43%  " }
    43 = @{ "B" = "Maker";   "C" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";   "D" = "2.759.39";  "E" = "  +5.07%  " }
    44 = @{                     "E" = "  +0.98%  " }
    45 = @{                     "E" = "  +5.98%  " }
    46 = @{ "B" = "Arweave"; "C" = "https://coinranking.com/coin/7XWg41D1+arweave-ar";       "D" = "35.58";     "E" = "  +26.15%  " }
    47 = @{ "B" = "USDe";    "C" = "https://coinranking.com/coin/exbfr2U-0+usde-usde";       "D" = "0.999";     "E" = "  +0.00%  " }
    48 = @{ "D" = "121.91";     "E" = "  +0.49%  " }
    49 = @{                     "E" = "  +2.02%  " }
    50 = @{                     "E" = "  +1.93%  " }
    51 = @{ "D" = "23.44";      "E" = "  +1.13%  " }
}

# Values such as "138.26" or "0.0991" round-trip as plain numbers through
# Range.Value, which would silently drop the source's fixed text
# formatting (trailing zeros, thousands-style dots, etc). Detect that case
# and force the cell to Text for the write, then restore the default
# "Normal" style so no stray number-format residue is left on the cell.
function Looks-Numeric($text) {
    return $text -match '^[+-]?\d+(\.\d+)*$'
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        $text = $rowData[$col]
        if (Looks-Numeric $text) {
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.Style = "Normal"
        } else {
            $cell.Value = $text
        }
    }
}
